{"js": "// EZ-3083 imail collection fixes\n//\n// The \"@Address5@\" placeholder paragraph is removed (merging it into the\n// following \"@Postcode@\" paragraph), and the document's \"_GoBack\"\n// bookmark is relocated to sit right before \"@Postcode@\" (where the\n// \"@Address5@\" run used to start).\n\n// 1. Drop the existing \"_GoBack\" bookmark wherever it currently lives \u2014\n//    it gets re-created at the new location below. (No-op if absent.)\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // Bookmark not present - nothing to remove.\n}\n\n// 2. Find the \"@Address5@\" paragraph and delete it outright; Word merges\n//    its content into the following paragraph (\"@Postcode@\") for us.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"@Address5@\") {\n    addressParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (addressParagraph) {\n  addressParagraph.delete();\n  await context.sync();\n}\n\n// 3. Re-load paragraphs and insert the \"_GoBack\" bookmark at the very\n//    start of the paragraph that now holds \"@Postcode@\".\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet postcodeParagraph = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"@Postcode@\") {\n    postcodeParagraph = paragraphs2.items[i];\n    break;\n  }\n}\n\nif (postcodeParagraph) {\n  const startRange = postcodeParagraph.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# EZ-3083 imail collection fixes\n#\n# The \"@Address5@\" placeholder paragraph is removed (its text plus the\n# paragraph mark that followed it are deleted, merging it into the\n# following \"@Postcode@\" paragraph), and the document's \"_GoBack\"\n# bookmark is relocated to sit right before \"@Postcode@\" (where the\n# \"@Address5@\" run used to start). Adding a bookmark named \"_GoBack\"\n# replaces/moves any existing \"_GoBack\" bookmark, so the stray one\n# elsewhere in the document goes away as a side effect of the Add call.\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"@Address5@\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $insertPos = $target.Range.Start\n    $nextPara = $target.Next()\n\n    # Delete from the start of the \"@Address5@\" paragraph through the\n    # start of the following paragraph - this removes the placeholder\n    # text together with its paragraph mark, merging the two paragraphs.\n    $mergeRange = $d.Range($insertPos, $nextPara.Range.Start)\n    $mergeRange.Delete()\n\n    # Re-create \"_GoBack\" collapsed at the merge point, right before\n    # \"@Postcode@\". This also removes/replaces any pre-existing\n    # \"_GoBack\" bookmark elsewhere in the document.\n    $bmRange = $d.Range($insertPos, $insertPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
